$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Rostin
$ws.Range("A4").Value = "Rostin"
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = "rostin@gmail.com"
$ws.Range("D4").Value = "Guarne"
$ws.Range("E4").Value = 31245456

# Row 5: Cristian Franco
$ws.Range("A5").Value = "Cristian Franco"
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = "cristianraigosa@gmail.com"
$ws.Range("D5").Value = "Medellín"
# Celular stored as text (large number), leading apostrophe forces text entry like in Excel
$ws.Range("E5").Value = "'3006487895"
